$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.046.68"
$ws.Range("E2").Value = "  -6.05%  "

$ws.Range("D3").Value = "2.449.57"
$ws.Range("E3").Value = "  -8.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.31"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.93%  "

$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("E8").Value = "  -3.59%  "

$ws.Range("D9").Value = "2.466.13"
$ws.Range("E9").Value = "  -8.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0992"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.158"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.35"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.67%  "

$ws.Range("D14").Value = "2.887.92"
$ws.Range("E14").Value = "  -8.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.04"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -9.24%  "

$ws.Range("D16").Value = "58.936.66"
$ws.Range("E16").Value = "  -6.07%  "

$ws.Range("E17").Value = "  -6.10%  "

$ws.Range("D18").Value = "2.516.99"
$ws.Range("E18").Value = "  -6.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.11"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -6.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.35"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.83"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.76%  "

$ws.Range("E22").Value = "  -3.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.71"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -9.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.459"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -9.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "60.71"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.26%  "

$ws.Range("E26").Value = "  -5.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.977"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.71"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.68%  "

$ws.Range("E29").Value = "  -9.89%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.82"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.04%  "

$ws.Range("E31").Value = "  -9.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.69"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -7.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.71"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.13%  "

$ws.Range("E35").Value = "  -7.89%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.44"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.47"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -9.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.71"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "317.49"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -10.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.87"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.20"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.837"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -11.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.71"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -7.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.71"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.77%  "

$ws.Range("E46").Value = "  -3.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.580"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.00%  "

$ws.Range("E48").Value = "  -6.20%  "

$ws.Range("E49").Value = "  -5.30%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.63"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.40%  "

$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.87"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -9.85%  "
